$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row appended to the bottom of the trade log (row 5)
$ws.Range("A5").Value = 42636.59270833333
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 10115.89
$ws.Range("D5").Value = 10085.129999999999
$ws.Range("E5").Value = 81.97
$ws.Range("F5").Value = 81.47
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = -0.61
$ws.Range("I5").Value = $false

# Match the date-formatted style used by the rest of the trade rows
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$excel.CutCopyMode = $false
